$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest GitHub Actions refresh.
# D-column price cells are textual (European dotted formatting like "28.510.87")
# so we force text format before assignment to avoid Excel auto-converting them
# to numbers, then reset the style back to Normal so no stray formatting sticks.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.510.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "

# Row 6
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5109"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.57%  "

# Row 8
$ws.Range("E8").Value = "  -2.24%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08216"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.58%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.92%  "

# Row 14
$ws.Range("E14").Value = "  +0.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.828.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.24%  "

# Row 17
$ws.Range("E17").Value = "  +3.91%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "

# Row 20
$ws.Range("E20").Value = "  +0.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.090"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.538.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "

# Row 24
$ws.Range("E24").Value = "  +2.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.17%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.032.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
$ws.Range("E29").Value = "  -2.09%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.64%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.120"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1093"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.763"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.654"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.37%  "

# Row 35
$ws.Range("E35").Value = "  -3.78%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2231"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.276"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.26%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02352"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.838"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6327"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.54%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.43%  "

# Row 42
$ws.Range("E42").Value = "  -0.65%  "

# Row 43
$ws.Range("E43").Value = "  +0.00%  "

# Row 44
$ws.Range("E44").Value = "  -0.14%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.35%  "

# Row 46
$ws.Range("E46").Value = "  +1.24%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.730"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.84%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.19%  "

# Row 49
$ws.Range("E49").Value = "  -0.77%  "

# Row 50
$ws.Range("E50").Value = "  -0.94%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06899"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
